$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7 (columns D, J, K, L, M, P), result of a row-data permutation
$values = @{
    2 = @{ D = 44839; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    3 = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    4 = @{ D = 44804; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    5 = @{ D = 44797; J = 60; K = 12000; L = 13000; M = 12500; P = 962 }
    6 = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    7 = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
}

foreach ($row in $values.Keys) {
    $rowData = $values[$row]
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("J$row").Value = $rowData.J
    $ws.Range("K$row").Value = $rowData.K
    $ws.Range("L$row").Value = $rowData.L
    $ws.Range("M$row").Value = $rowData.M
    $ws.Range("P$row").Value = $rowData.P
}
